$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as pure text, avoiding Excel auto-numeric-conversion
# for values that look like numbers (e.g. "1.001", "44.70", "0.00001081"),
# while leaving the cell style index untouched (no "s" attribute added).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '22.301.33'
$ws.Range('E2').Value = '  -1.18%  '

# Row 3
Set-TextValue 'D3' '1.559.52'
$ws.Range('E3').Value = '  -1.20%  '

# Row 4
$ws.Range('E4').Value = '  -0.14%  '

# Row 5
Set-TextValue 'D5' '1.001'
$ws.Range('E5').Value = '  -0.08%  '

# Row 6
$ws.Range('E6').Value = '  +0.19%  '

# Row 7
Set-TextValue 'D7' '0.3796'
$ws.Range('E7').Value = '  +2.34%  '

# Row 8
Set-TextValue 'D8' '0.3291'
$ws.Range('E8').Value = '  -1.83%  '

# Row 9
Set-TextValue 'D9' '44.90'
$ws.Range('E9').Value = '  -7.66%  '

# Row 10
Set-TextValue 'D10' '1.143'
$ws.Range('E10').Value = '  -0.18%  '

# Row 11
Set-TextValue 'D11' '0.07398'
$ws.Range('E11').Value = '  -1.34%  '

# Row 12
$ws.Range('E12').Value = '  -0.11%  '

# Row 13
Set-TextValue 'D13' '20.36'
$ws.Range('E13').Value = '  -3.34%  '

# Row 14
Set-TextValue 'D14' '5.880'
$ws.Range('E14').Value = '  -2.08%  '

# Row 15
$ws.Range('E15').Value = '  -2.57%  '

# Row 16
Set-TextValue 'D16' '1.553.16'
$ws.Range('E16').Value = '  -1.81%  '

# Row 17
Set-TextValue 'D17' '0.00001081'
$ws.Range('E17').Value = '  -3.82%  '

# Row 18
Set-TextValue 'D18' '0.06661'
$ws.Range('E18').Value = '  -1.77%  '

# Row 19
Set-TextValue 'D19' '86.59'
$ws.Range('E19').Value = '  -2.37%  '

# Row 20
Set-TextValue 'D20' '6.443'
$ws.Range('E20').Value = '  +0.32%  '

# Row 21
Set-TextValue 'D21' '1.001'
$ws.Range('E21').Value = '  -0.08%  '

# Row 22
Set-TextValue 'D22' '16.25'
$ws.Range('E22').Value = '  -2.12%  '

# Row 23
Set-TextValue 'D23' '11.78'
$ws.Range('E23').Value = '  -3.28%  '

# Row 24
Set-TextValue 'D24' '22.283.46'
$ws.Range('E24').Value = '  -1.19%  '

# Row 25
$ws.Range('E25').Value = '  -4.48%  '

# Row 26
Set-TextValue 'D26' '2.581'
$ws.Range('E26').Value = '  -0.80%  '

# Row 27
Set-TextValue 'D27' '151.20'
$ws.Range('E27').Value = '  -1.13%  '

# Row 28
$ws.Range('E28').Value = '  -2.02%  '

# Row 29
Set-TextValue 'D29' '4.944'
$ws.Range('E29').Value = '  -1.67%  '

# Row 30
Set-TextValue 'D30' '123.33'
$ws.Range('E30').Value = '  -0.91%  '

# Row 31
Set-TextValue 'D31' '1.738.40'
$ws.Range('E31').Value = '  -1.03%  '

# Row 32
Set-TextValue 'D32' '1.084'
$ws.Range('E32').Value = '  +2.07%  '

# Row 33
Set-TextValue 'D33' '5.952'
$ws.Range('E33').Value = '  -4.10%  '

# Row 34
Set-TextValue 'D34' '1.927'
$ws.Range('E34').Value = '  -4.36%  '

# Row 35
Set-TextValue 'D35' '9.479'
$ws.Range('E35').Value = '  -2.37%  '

# Row 36
Set-TextValue 'D36' '0.08221'
$ws.Range('E36').Value = '  -1.32%  '

# Row 37
Set-TextValue 'D37' '0.02368'
$ws.Range('E37').Value = '  -4.09%  '

# Row 38
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D38' '5.400'
$ws.Range('E38').Value = '  -0.64%  '

# Row 39
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D39' '0.06376'
$ws.Range('E39').Value = '  -0.44%  '

# Row 40
Set-TextValue 'D40' '0.2170'
$ws.Range('E40').Value = '  -5.57%  '

# Row 41
Set-TextValue 'D41' '1.247'
$ws.Range('E41').Value = '  -3.97%  '

# Row 42
$ws.Range('E42').Value = '  -2.80%  '

# Row 43
Set-TextValue 'D43' '0.6101'
$ws.Range('E43').Value = '  -4.13%  '

# Row 44
Set-TextValue 'D44' '1.000'
$ws.Range('E44').Value = '  -0.15%  '

# Row 45
Set-TextValue 'D45' '13.83'
$ws.Range('E45').Value = '  -1.20%  '

# Row 46
Set-TextValue 'D46' '3.765'
$ws.Range('E46').Value = '  -0.30%  '

# Row 47
Set-TextValue 'D47' '0.5922'
$ws.Range('E47').Value = '  -4.96%  '

# Row 48
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D48' '123.75'
$ws.Range('E48').Value = '  -1.17%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '1.988'
$ws.Range('E49').Value = '  -3.91%  '

# Row 50
Set-TextValue 'D50' '1.182'
$ws.Range('E50').Value = '  -3.34%  '

# Row 51
Set-TextValue 'D51' '0.07090'
$ws.Range('E51').Value = '  -2.71%  '
